# Weekly update: insert the newest "Macroferia Regional de Talca - Zanahoria"
# market record (2021-10-20, Región de Ñuble) at the top of the date-ordered
# block, pushing the existing rows 138-201 down to 139-202.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 138 (shifts 138:201 -> 139:202).
$ws.Rows.Item(138).Insert()

# Populate the new row with the latest weekly observation.
$ws.Range("A138").Value2 = 5
$ws.Range("B138").Value2 = "Macroferia Regional de Talca"
$ws.Range("C138").Value2 = "Maule"
$ws.Range("D138").Value2 = 44489
$ws.Range("E138").Value2 = 7
$ws.Range("F138").Value2 = 100114013
$ws.Range("G138").Value2 = "Zanahoria"
$ws.Range("H138").Value2 = "Sin especificar"
$ws.Range("I138").Value2 = "Primera"
$ws.Range("J138").Value2 = 300
$ws.Range("K138").Value2 = 8000
$ws.Range("L138").Value2 = 8000
$ws.Range("M138").Value2 = 8000
$ws.Range("N138").Value2 = "$/saco 20 kilos"
$ws.Range("O138").Value2 = "Región de Ñuble"
$ws.Range("P138").Value2 = 400
$ws.Range("Q138").Value2 = 20
$ws.Range("R138").Value2 = "Hortaliza"
